# Update the build timestamp embedded in the version strings throughout the
# workbook, from "February 03 2026 17.29.55 EST" to
# "February 03 2026 18.05.36 EST".

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: version banner (A2) and citation text (A6)
$a2 = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# "Boundaries and methane sources" sheet: build_version column (S2:S7)
for ($row = 2; $row -le 7; $row++) {
    $cell = $wsData.Range("S" + $row)
    $current = $cell.Value()
    $cell.Value = $current.Replace($oldStamp, $newStamp)
}
